$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N; the old N/O/P columns (Late /
# heading / Outstanding) shift right to O/P/Q.
$ws.Columns("N").Insert()

# The newly inserted column picks up the width of its left neighbour
# (column M), matching how Excel renders a freshly inserted column.
$ws.Columns("N").ColumnWidth = 10.17

# Make "Repayment schedule" the active sheet/tab and select cell S8 on it,
# matching the saved workbook view state (was "Transactions" before).
$ws.Activate()
$ws.Range("S8").Select() | Out-Null
